# Update the Sheet1 worksheet to match the target revision:
#  - header cell A1 text "product_no" -> "material_code"
#  - resize columns A-D (E stays the same)
#  - move the active selection from B1 to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Rename header cell A1
$ws.Range("A1").Value = "material_code"

# 2) Resize columns A:D to their new widths (column E is unchanged)
#    The values below are chosen so that, after the runtime's internal
#    pixel-rounding of ColumnWidth, the stored widths land as close as
#    possible to the target widths (104/7, 54/7, 100/7 and 6 characters).
$ws.Columns.Item(1).ColumnWidth = 14.0
$ws.Columns.Item(2).ColumnWidth = 6.833333333333334
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 5.166666666666666

# 3) Move the selection to A3
$ws.Range("A3").Select()
